$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.908.69"
$ws.Range("E2").Value = "  -5.24%  "

$ws.Range("D3").Value = "3.309.10"
$ws.Range("E3").Value = "  -5.34%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.01"
$ws.Range("E5").Value = "  -3.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.21"
$ws.Range("E6").Value = "  -5.55%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "3.309.97"
$ws.Range("E8").Value = "  -5.31%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.475"
$ws.Range("E9").Value = "  -2.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.24"
$ws.Range("E10").Value = "  -4.97%  "

$ws.Range("E11").Value = "  -6.22%  "

$ws.Range("E12").Value = "  -4.03%  "

$ws.Range("D13").Value = "3.880.72"
$ws.Range("E13").Value = "  -5.17%  "

$ws.Range("E14").Value = "  -1.55%  "

$ws.Range("D15").Value = "3.318.62"
$ws.Range("E15").Value = "  -5.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000166"
$ws.Range("E16").Value = "  -7.71%  "

$ws.Range("D17").Value = "61.030.61"
$ws.Range("E17").Value = "  -5.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "24.22"
$ws.Range("E18").Value = "  -4.18%  "

$ws.Range("E19").Value = "  -3.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.02"
$ws.Range("E20").Value = "  -10.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.13"
$ws.Range("E21").Value = "  -2.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "347.64"
$ws.Range("E22").Value = "  -10.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.550"
$ws.Range("E23").Value = "  -5.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("E24").Value = "  -0.28%  "

$ws.Range("D25").Value = "3.445.54"
$ws.Range("E25").Value = "  -5.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.80"
$ws.Range("E26").Value = "  -5.83%  "

$ws.Range("E27").Value = "  -8.30%  "

$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.09"
$ws.Range("E29").Value = "  -3.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.43"
$ws.Range("E30").Value = "  -4.28%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.78"
$ws.Range("E31").Value = "  -4.53%  "

$ws.Range("E32").Value = "  -6.81%  "

$ws.Range("E33").Value = "  -0.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.146"
$ws.Range("E34").Value = "  -5.49%  "

$ws.Range("D35").Value = "3.343.38"
$ws.Range("E35").Value = "  -5.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.26"
$ws.Range("E36").Value = "  -4.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.35"
$ws.Range("E37").Value = "  +0.40%  "

$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.54"
$ws.Range("E38").Value = "  -1.77%  "

$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.62"
$ws.Range("E39").Value = "  -4.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.46"
$ws.Range("E40").Value = "  -4.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0747"
$ws.Range("E41").Value = "  -4.83%  "

$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.85"
$ws.Range("E43").Value = "  -2.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.742"
$ws.Range("E44").Value = "  -8.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.19"
$ws.Range("E45").Value = "  -5.02%  "

$ws.Range("E46").Value = "  -6.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.53"
$ws.Range("E47").Value = "  -6.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.42"
$ws.Range("E48").Value = "  -8.13%  "

$ws.Range("E49").Value = "  -3.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.844"
$ws.Range("E50").Value = "  -7.93%  "

$ws.Range("D51").Value = "2.190.00"
$ws.Range("E51").Value = "  -9.66%  "

